$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("urls")

$ws.Range("C2").Value = 4869.300000000134
$ws.Range("E2").Value = 16195

$ws.Range("C3").Value = 4768.7999999999
$ws.Range("D3").Value = 0.000723548
$ws.Range("E3").Value = 18477

$ws.Range("C4").Value = 4735.900000000067
$ws.Range("D4").Value = 0.0030348684
$ws.Range("E4").Value = 18292.333333333332

$ws.Range("C5").Value = 4735.300000000033
$ws.Range("D5").Value = 0.0030348684
$ws.Range("E5").Value = 18277.333333333332

$ws.Range("C6").Value = 4636.566666666734
$ws.Range("D6").Value = 0.0022970415
$ws.Range("E6").Value = 19603.333333333332

$ws.Range("C7").Value = 4635.366666666533
$ws.Range("D7").Value = 0.0022970415000000003
$ws.Range("E7").Value = 19502.333333333332

$ws.Range("C8").Value = 4535.5333333332665
$ws.Range("E8").Value = 15587.666666666666

$ws.Range("C9").Value = 4602.7666666666
$ws.Range("E9").Value = 15527.666666666666
